$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B3"  = 4150.640883136433
    "B4"  = 9676.511341654197
    "B5"  = 228381.1975844747
    "B6"  = 5496.795799658239
    "B7"  = 4150.64088313644
    "B8"  = 9676.511341654201
    "B9"  = 228381.1975831594
    "B10" = 5496.795799658239
    "B15" = 4067.628065473661
    "B16" = 9579.746228237562
    "B17" = 105.1082831648694
    "B18" = 83.01281766272777
    "B19" = 96.76511341654104
    "B20" = 3.25077164427432
    "B21" = 228381.1975808766
    "B22" = 4150.640883136389
    "B23" = 9676.511341654104
    "B24" = 228381.1975808766
    "B25" = 108.3590548091438
    "B26" = 242208.349805667
    "B27" = 4876.15746641147
    "B28" = 27369.54352804038
    "B29" = 6021.299576168884
    "B30" = 4026.951784818923
    "B31" = 95.79746228237563
    "B32" = 40.67628065473661
    "B33" = 9483.948765955189
    "B34" = 4067.62806547366
    "B35" = 9579.746228237564
    "B36" = 13647.37429371122
    "B37" = 10508.47820615764
    "B38" = 2311.865205354681
    "B39" = 4026.951784818923
    "B40" = 9483.948765955189
    "B41" = 5496.795799658239
    "B42" = 228564.2262836001
    "B43" = 738.1143148937871
    "B44" = 25607.45361530723
    "B46" = 13211.96224793503
    "B47" = 24990.07656282882
    "B51" = 617.3770524784122
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
